$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds issue IDs that look numeric ("17"). Format the cell as text
# first so Excel stores the new value as a string (matching every other row
# in the sheet) instead of auto-converting it to a number.
$ws.Range("A13").NumberFormat = "@"

$ws.Range("A13").Value = "17"
$ws.Range("B13").Value = "FR_OPERATIONS issue"
$ws.Range("C13").Value = "open"
$ws.Range("D13").Value = "2025-03-25T11:47:03Z"
$ws.Range("E13").Value = "bug"
